$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")
$fmtSrc = $wb.Worksheets.Item("AvailableIndices")

# --- Insert a new top row; this shifts the existing header + 2 data rows down to rows 2-4 ---
$ws.Rows.Item(1).Insert() | Out-Null

# Row 1: informational note, column D only
$ws.Range("D1").Value = "AKShare did not support all HKEX data."

# --- Header row (now row 2): existing "akshare" header plus new "A股"/"港股" columns ---
$ws.Range("D2").Value = "akshare"
$ws.Range("E2").Value = "A股"
$ws.Range("F2").Value = "港股"

# --- Existing data rows (now rows 3 and 4): both are A-share ("x" in col E) only ---
$ws.Range("D3").Value = "AKShareAvailableIndicesData"
$ws.Range("E3").Value = "x"

$ws.Range("D4").Value = "AKShareBalanceSheetData"
$ws.Range("E4").Value = "x"

# --- New rows 5-8 ---
$ws.Range("A5").Value = "CompanyNewsData"
$ws.Range("D5").Value = "AKShareCompanyNewsData"
$ws.Range("E5").Value = "x"

$ws.Range("A6").Value = "EquityHistoricalData"
$ws.Range("D6").Value = "AKShareEquityHistoricalData"
$ws.Range("E6").Value = "x"
$ws.Range("F6").Value = "x"

$ws.Range("A7").Value = "EquityQuoteData"
$ws.Range("D7").Value = "AKShareEquityQuoteData"
$ws.Range("E7").Value = "x"
$ws.Range("F7").Value = "x"

$ws.Range("A8").Value = "HistoricalDividendsData"
$ws.Range("D8").Value = "AKShareHistoricalDividendsData"
$ws.Range("E8").Value = "x"
$ws.Range("F8").Value = "x"

# --- Formatting: match the look of the header/body styles used elsewhere in the workbook ---
# Bold + thin border on the header row (D2:F2), same as AvailableIndices!A1
$fmtSrc.Range("A1").Copy() | Out-Null
$ws.Range("D2:F2").PasteSpecial(-4122) | Out-Null

# Thin border (no bold) on every body cell in columns D:F, same as AvailableIndices!A2
$fmtSrc.Range("A2").Copy() | Out-Null
$ws.Range("D3:F8").PasteSpecial(-4122) | Out-Null

# --- Selection / active sheet bookkeeping, matching the saved view state ---
$ws.Activate() | Out-Null
$ws.Range("E3").Select() | Out-Null
